# Apply the "baptismal certificate" edit described in the commit:
#  - sheet1 ("Confirmation Events"): widen column A, add a new event row
#    (row 6) for "Upload Baptismal Certificate"
#  - sheet2 ("Candidates with events"): add two new event slots
#    (candidate_events.3 / candidate_events.4), i.e. 4 new columns T:W
#    mirroring the existing N:O / P:Q / R:S pattern

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Confirmation Events"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Widen column A to fit the new, longer event name.
$ws1.Columns.Item(1).ColumnWidth = 33.86

# Row 6 currently holds blank, formatted placeholder cells (A6:E6).
# Copy the formatting from row 5 (the previous event row) down into row 6
# so the new row picks up the same styles (text / date / text), then fill
# in the values.
$ws1.Range("A5:C5").Copy($ws1.Range("A6:C6"))

$ws1.Cells.Item(6, 1).Value = "Upload Baptismal Certificate"
$ws1.Cells.Item(6, 2).Value2 = 42594
$ws1.Cells.Item(6, 3).Value = "<p><em><strong>Upload certificate</strong></em></p>"

# ---------------------------------------------------------------------
# Sheet 2: "Candidates with events"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# New columns T:W get the same display width as the existing event
# columns (Q:S, 18.6719 characters).
$ws2.Columns.Item(20).ColumnWidth = 17.93
$ws2.Columns.Item(21).ColumnWidth = 17.93
$ws2.Columns.Item(22).ColumnWidth = 17.93
$ws2.Columns.Item(23).ColumnWidth = 17.93

# Header row (row 1): T1:W1 mirror the text-header style used in N1:S1.
$ws2.Range("N1").Copy($ws2.Range("T1"))
$ws2.Range("N1").Copy($ws2.Range("U1"))
$ws2.Range("N1").Copy($ws2.Range("V1"))
$ws2.Range("N1").Copy($ws2.Range("W1"))

$ws2.Cells.Item(1, 20).Value = "candidate_events.3.completed_date"
$ws2.Cells.Item(1, 21).Value = "candidate_events.3.verified"
$ws2.Cells.Item(1, 22).Value = "candidate_events.4.completed_date"
$ws2.Cells.Item(1, 23).Value = "candidate_events.4.verified"

# Row 2 (Vicki Kristoff): mirror R2 (plain/empty) and S2 (boolean False)
# styles for the two new completed_date/verified pairs.
$ws2.Range("R2").Copy($ws2.Range("T2"))
$ws2.Range("S2").Copy($ws2.Range("U2"))
$ws2.Cells.Item(2, 21).Value = $false
$ws2.Range("R2").Copy($ws2.Range("V2"))
$ws2.Range("S2").Copy($ws2.Range("W2"))
$ws2.Cells.Item(2, 23).Value = $false

# Row 3 (Paul Kristoff): mirror P3 (empty date-formatted) and Q3 (boolean
# False) styles for the two new completed_date/verified pairs.
$ws2.Range("P3").Copy($ws2.Range("T3"))
$ws2.Cells.Item(3, 20).ClearContents()
$ws2.Range("Q3").Copy($ws2.Range("U3"))
$ws2.Cells.Item(3, 21).Value = $false
$ws2.Range("P3").Copy($ws2.Range("V3"))
$ws2.Cells.Item(3, 22).ClearContents()
$ws2.Range("Q3").Copy($ws2.Range("W3"))
$ws2.Cells.Item(3, 23).Value = $false

# Row 4 (foobar) and the remaining blank rows 5:10: mirror the plain
# empty style used throughout columns N:S on those rows.
$ws2.Range("N4:Q4").Copy($ws2.Range("T4:W4"))
$ws2.Range("N5:Q5").Copy($ws2.Range("T5:W5"))
$ws2.Range("N6:Q6").Copy($ws2.Range("T6:W6"))
$ws2.Range("N7:Q7").Copy($ws2.Range("T7:W7"))
$ws2.Range("N8:Q8").Copy($ws2.Range("T8:W8"))
$ws2.Range("N9:Q9").Copy($ws2.Range("T9:W9"))
$ws2.Range("N10:Q10").Copy($ws2.Range("T10:W10"))
